$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.329.24'
$ws.Range("E2").Value = '  +7.72%  '
$ws.Range("D3").Value = '3.403.77'
$ws.Range("E3").Value = '  +4.94%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '411.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '121.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.94%  '
$ws.Range("D7").Value = '3.398.25'
$ws.Range("E7").Value = '  +4.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.636'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.122'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +28.93%  '
$ws.Range("E12").Value = '  +4.78%  '
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '3.937.49'
$ws.Range("E14").Value = '  +4.86%  '
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").Value = '3.419.84'
$ws.Range("E17").Value = '  +5.02%  '
$ws.Range("D18").Value = '61.267.69'
$ws.Range("E18").Value = '  +8.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("E21").Value = '  +12.14%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '297.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '76.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("E27").Value = '  +9.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.89%  '
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.171'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("E32").Value = '  +5.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +18.74%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0476'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("E42").Value = '  +5.06%  '
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '133.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.281'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.05%  '
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("D50").Value = '2.201.25'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("D51").Value = '3.742.55'
$ws.Range("E51").Value = '  +5.03%  '
